# Add the new "2023 - Høst" semester row beneath the existing "2023 - Vår" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Value = "2023 - Høst"
$ws.Range("B15").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/hjemme-23-h.pdf)"
$ws.Range("C15").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-23-h-solprop.html)"
$ws.Range("D15").Value = "[Materiale](tidligere-eksamensoppgaver/hjemme-23-h-ekstra.zip)"

# Match the saved cursor/selection position recorded in the committed file.
$ws.Range("D19").Select() | Out-Null
